# Edit script generated from the target diff.
# Applies cell-level changes to worksheet "Artfynd" (the active/only sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Part 1: existing rows whose field values were re-shuffled / corrected.
# ----------------------------------------------------------------------
# --- Row 5 ---
$ws.Range("A5").Value = 131260583
$ws.Range("M5").Value = "färska spår"
$ws.Range("Q5").Value = 488834
$ws.Range("R5").Value = 6665228
$ws.Range("Z5").Value = "15:30"
$ws.Range("AB5").Value = "15:30"
$ws.Range("AC5").Value = "Ringhack på tall."
# --- Row 6 ---
$ws.Range("A6").Value = 131256691
$ws.Range("M6").Value = "äldre spår"
$ws.Range("Q6").Value = 488667
$ws.Range("R6").Value = 6665262
$ws.Range("Z6").Value = "10:55"
$ws.Range("AB6").Value = "10:55"
$ws.Range("AC6").Value = "Ringhack på gran."
# --- Row 12 ---
$ws.Range("A12").Value = 131257520
$ws.Range("B12").Value = 79244
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = "Garnlav"
$ws.Range("G12").Value = "Alectoria sarmentosa"
$ws.Range("H12").Value = "(Ach.) Ach."
$ws.Range("M12").ClearContents()
$ws.Range("Q12").Value = 488939
$ws.Range("R12").Value = 6665149
$ws.Range("Z12").Value = "11:41"
$ws.Range("AB12").Value = "11:41"
$ws.Range("AC12").Value = "Gran"
# --- Row 13 ---
$ws.Range("A13").Value = 131256673
$ws.Range("B13").Value = 57884
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = "Tretåig hackspett"
$ws.Range("G13").Value = "Picoides tridactylus"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("M13").Value = "äldre spår"
$ws.Range("Q13").Value = 488652
$ws.Range("R13").Value = 6665282
$ws.Range("Z13").Value = "10:54"
$ws.Range("AB13").Value = "10:54"
$ws.Range("AC13").Value = "Ringhack på tall."
# --- Row 14 ---
$ws.Range("A14").Value = 131260641
$ws.Range("Q14").Value = 488859
$ws.Range("R14").Value = 6665292
$ws.Range("Z14").Value = "15:34"
$ws.Range("AB14").Value = "15:34"
$ws.Range("AC14").Value = "Ringhack på gran."
# --- Row 15 ---
$ws.Range("A15").Value = 131257290
$ws.Range("Q15").Value = 488842
$ws.Range("R15").Value = 6665224
$ws.Range("Z15").Value = "11:26"
$ws.Range("AB15").Value = "11:26"
# --- Row 25 ---
$ws.Range("A25").Value = 131257650
$ws.Range("Q25").Value = 488911
$ws.Range("R25").Value = 6665227
$ws.Range("Z25").Value = "12:00"
$ws.Range("AB25").Value = "12:00"
$ws.Range("AC25").Value = "Gran"
# --- Row 26 ---
$ws.Range("A26").Value = 131257045
$ws.Range("Q26").Value = 488760
$ws.Range("R26").Value = 6665302
$ws.Range("Z26").Value = "11:10"
$ws.Range("AB26").Value = "11:10"
$ws.Range("AC26").Value = "Gran."
# --- Row 28 ---
$ws.Range("A28").Value = 131257239
$ws.Range("B28").Value = 57884
$ws.Range("E28").Value = 100109
$ws.Range("F28").Value = "Tretåig hackspett"
$ws.Range("G28").Value = "Picoides tridactylus"
$ws.Range("H28").Value = "(Linnaeus, 1758)"
$ws.Range("M28").Value = "färska spår"
$ws.Range("Q28").Value = 488852
$ws.Range("R28").Value = 6665286
$ws.Range("Z28").Value = "11:23"
$ws.Range("AB28").Value = "11:23"
$ws.Range("AC28").Value = "Barkfläk, hagelsalva."
# --- Row 29 ---
$ws.Range("A29").Value = 131256649
$ws.Range("Q29").Value = 488685
$ws.Range("R29").Value = 6665288
$ws.Range("Z29").Value = "10:52"
$ws.Range("AB29").Value = "10:52"
$ws.Range("AC29").Value = "Gran."
# --- Row 30 ---
$ws.Range("A30").Value = 131255910
$ws.Range("B30").Value = 79244
$ws.Range("E30").Value = 6425
$ws.Range("F30").Value = "Garnlav"
$ws.Range("G30").Value = "Alectoria sarmentosa"
$ws.Range("H30").Value = "(Ach.) Ach."
$ws.Range("M30").ClearContents()
$ws.Range("Q30").Value = 488763
$ws.Range("R30").Value = 6665157
$ws.Range("Z30").Value = "10:03"
$ws.Range("AB30").Value = "10:03"
$ws.Range("AC30").Value = "Tall."
# --- Row 32 ---
$ws.Range("A32").Value = 131257659
$ws.Range("B32").Value = 57884
$ws.Range("E32").Value = 100109
$ws.Range("F32").Value = "Tretåig hackspett"
$ws.Range("G32").Value = "Picoides tridactylus"
$ws.Range("H32").Value = "(Linnaeus, 1758)"
$ws.Range("M32").Value = "äldre spår"
$ws.Range("Q32").Value = 488901
$ws.Range("R32").Value = 6665222
$ws.Range("Z32").Value = "12:02"
$ws.Range("AB32").Value = "12:02"
$ws.Range("AC32").Value = "Ringhack på tall."
# --- Row 34 ---
$ws.Range("A34").Value = 131258537
$ws.Range("B34").Value = 79244
$ws.Range("E34").Value = 6425
$ws.Range("F34").Value = "Garnlav"
$ws.Range("G34").Value = "Alectoria sarmentosa"
$ws.Range("H34").Value = "(Ach.) Ach."
$ws.Range("M34").ClearContents()
$ws.Range("Q34").Value = 488726
$ws.Range("R34").Value = 6665209
$ws.Range("Z34").Value = "13:02"
$ws.Range("AB34").Value = "13:02"
$ws.Range("AC34").Value = "Gran"
# --- Row 41 ---
$ws.Range("A41").Value = 131257343
$ws.Range("B41").Value = 57884
$ws.Range("E41").Value = 100109
$ws.Range("F41").Value = "Tretåig hackspett"
$ws.Range("G41").Value = "Picoides tridactylus"
$ws.Range("H41").Value = "(Linnaeus, 1758)"
$ws.Range("M41").Value = "färska spår"
$ws.Range("Q41").Value = 488872
$ws.Range("R41").Value = 6665191
$ws.Range("Z41").Value = "11:29"
$ws.Range("AB41").Value = "11:29"
$ws.Range("AC41").Value = "Mycket barkfläk, hagelsalvor på många träd, skalade klena senvuxna granar."
# --- Row 42 ---
$ws.Range("A42").Value = 131257035
$ws.Range("B42").Value = 79244
$ws.Range("E42").Value = 6425
$ws.Range("F42").Value = "Garnlav"
$ws.Range("G42").Value = "Alectoria sarmentosa"
$ws.Range("H42").Value = "(Ach.) Ach."
$ws.Range("M42").ClearContents()
$ws.Range("Q42").Value = 488760
$ws.Range("R42").Value = 6665301
$ws.Range("Z42").Value = "11:09"
$ws.Range("AB42").Value = "11:09"
$ws.Range("AC42").Value = "Gran"

# ----------------------------------------------------------------------
# Part 2: three new observation rows appended at the bottom (43-45).
# ----------------------------------------------------------------------
# --- Row 43 (new) ---
$ws.Range("A43").Value = 131273946
$ws.Range("B43").Value = 79244
$ws.Range("D43").Value = "NT"
$ws.Range("E43").Value = 6425
$ws.Range("F43").Value = "Garnlav"
$ws.Range("G43").Value = "Alectoria sarmentosa"
$ws.Range("H43").Value = "(Ach.) Ach."
$ws.Range("P43").Value = "Hyttfallet, Dlr"
$ws.Range("Q43").Value = 488774
$ws.Range("R43").Value = 6665353
$ws.Range("S43").Value = 50
$ws.Range("T43").Value = "Dalarna"
$ws.Range("U43").Value = "Ludvika"
$ws.Range("V43").Value = "Dalarna"
$ws.Range("W43").Value = "Grangärde"
$ws.Range("Y43").NumberFormat = "@"
$ws.Range("Y43").Value = "2026-02-22"
$ws.Range("AA43").NumberFormat = "@"
$ws.Range("AA43").Value = "2026-02-22"
$ws.Range("AD43").Value = $false
$ws.Range("AE43").Value = $false
$ws.Range("AG43").Value = $false
$ws.Range("AW43").Value = "Anna-Lena Thommson"
$ws.Range("AX43").Value = "Anna-Lena Thommson"

# --- Row 44 (new) ---
$ws.Range("A44").Value = 131273991
$ws.Range("B44").Value = 79244
$ws.Range("D44").Value = "NT"
$ws.Range("E44").Value = 6425
$ws.Range("F44").Value = "Garnlav"
$ws.Range("G44").Value = "Alectoria sarmentosa"
$ws.Range("H44").Value = "(Ach.) Ach."
$ws.Range("P44").Value = "Hyttfallet, Dlr"
$ws.Range("Q44").Value = 488928
$ws.Range("R44").Value = 6665146
$ws.Range("S44").Value = 50
$ws.Range("T44").Value = "Dalarna"
$ws.Range("U44").Value = "Ludvika"
$ws.Range("V44").Value = "Dalarna"
$ws.Range("W44").Value = "Grangärde"
$ws.Range("Y44").NumberFormat = "@"
$ws.Range("Y44").Value = "2026-02-22"
$ws.Range("AA44").NumberFormat = "@"
$ws.Range("AA44").Value = "2026-02-22"
$ws.Range("AD44").Value = $false
$ws.Range("AE44").Value = $false
$ws.Range("AG44").Value = $false
$ws.Range("AW44").Value = "Anna-Lena Thommson"
$ws.Range("AX44").Value = "Anna-Lena Thommson"

# --- Row 45 (new) ---
$ws.Range("A45").Value = 131273940
$ws.Range("B45").Value = 79244
$ws.Range("D45").Value = "NT"
$ws.Range("E45").Value = 6425
$ws.Range("F45").Value = "Garnlav"
$ws.Range("G45").Value = "Alectoria sarmentosa"
$ws.Range("H45").Value = "(Ach.) Ach."
$ws.Range("P45").Value = "Hyttfallet, Dlr"
$ws.Range("Q45").Value = 488785
$ws.Range("R45").Value = 6665325
$ws.Range("S45").Value = 50
$ws.Range("T45").Value = "Dalarna"
$ws.Range("U45").Value = "Ludvika"
$ws.Range("V45").Value = "Dalarna"
$ws.Range("W45").Value = "Grangärde"
$ws.Range("Y45").NumberFormat = "@"
$ws.Range("Y45").Value = "2026-02-22"
$ws.Range("AA45").NumberFormat = "@"
$ws.Range("AA45").Value = "2026-02-22"
$ws.Range("AD45").Value = $false
$ws.Range("AE45").Value = $false
$ws.Range("AG45").Value = $false
$ws.Range("AW45").Value = "Anna-Lena Thommson"
$ws.Range("AX45").Value = "Anna-Lena Thommson"

